$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.887
$ws.Range("B3").Value = 0.65
$ws.Range("B4").Value = 0.656

$ws.Rows.Item(5).Delete()
